{"js": "// Update the division-problem worksheet numbers.\n// Each entry is an exact \"old text\" -> \"new text\" replacement for a\n// <w:t> run inside the table cells (format \"NNN\u00f7D=\").\nconst replacements = [\n  [\"107\u00f78=\", \"323\u00f74=\"],\n  [\"808\u00f75=\", \"747\u00f75=\"],\n  [\"377\u00f77=\", \"678\u00f76=\"],\n  [\"529\u00f77=\", \"571\u00f74=\"],\n  [\"298\u00f77=\", \"706\u00f75=\"],\n  [\"420\u00f74=\", \"238\u00f76=\"],\n  [\"224\u00f76=\", \"369\u00f73=\"],\n  [\"321\u00f75=\", \"528\u00f73=\"],\n  [\"888\u00f76=\", \"375\u00f76=\"],\n  [\"836\u00f79=\", \"442\u00f74=\"],\n  [\"364\u00f75=\", \"557\u00f72=\"],\n  [\"139\u00f79=\", \"892\u00f75=\"],\n  [\"342\u00f72=\", \"444\u00f74=\"],\n  [\"639\u00f74=\", \"204\u00f76=\"],\n  [\"602\u00f76=\", \"157\u00f79=\"],\n  [\"113\u00f79=\", \"469\u00f77=\"],\n  [\"577\u00f74=\", \"896\u00f78=\"],\n  [\"888\u00f79=\", \"783\u00f72=\"],\n  [\"839\u00f78=\", \"182\u00f73=\"],\n  [\"568\u00f75=\", \"720\u00f76=\"],\n  [\"779\u00f77=\", \"212\u00f79=\"],\n  [\"230\u00f75=\", \"837\u00f72=\"],\n  [\"356\u00f76=\", \"561\u00f73=\"],\n  [\"771\u00f72=\", \"949\u00f77=\"],\n  [\"537\u00f77=\", \"953\u00f79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division-problem worksheet numbers.\n# Each entry is an exact \"old text\" -> \"new text\" replacement for the\n# text of a table cell (format \"NNN\u00f7D=\").\n\n# Word \"Find\" constants aren't predefined as PS variables in this host,\n# so define them explicitly.\n$wdReplaceNone = 0\n$wdReplaceOne = 1\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = @(\n    @{ Old = \"107\u00f78=\"; New = \"323\u00f74=\" },\n    @{ Old = \"808\u00f75=\"; New = \"747\u00f75=\" },\n    @{ Old = \"377\u00f77=\"; New = \"678\u00f76=\" },\n    @{ Old = \"529\u00f77=\"; New = \"571\u00f74=\" },\n    @{ Old = \"298\u00f77=\"; New = \"706\u00f75=\" },\n    @{ Old = \"420\u00f74=\"; New = \"238\u00f76=\" },\n    @{ Old = \"224\u00f76=\"; New = \"369\u00f73=\" },\n    @{ Old = \"321\u00f75=\"; New = \"528\u00f73=\" },\n    @{ Old = \"888\u00f76=\"; New = \"375\u00f76=\" },\n    @{ Old = \"836\u00f79=\"; New = \"442\u00f74=\" },\n    @{ Old = \"364\u00f75=\"; New = \"557\u00f72=\" },\n    @{ Old = \"139\u00f79=\"; New = \"892\u00f75=\" },\n    @{ Old = \"342\u00f72=\"; New = \"444\u00f74=\" },\n    @{ Old = \"639\u00f74=\"; New = \"204\u00f76=\" },\n    @{ Old = \"602\u00f76=\"; New = \"157\u00f79=\" },\n    @{ Old = \"113\u00f79=\"; New = \"469\u00f77=\" },\n    @{ Old = \"577\u00f74=\"; New = \"896\u00f78=\" },\n    @{ Old = \"888\u00f79=\"; New = \"783\u00f72=\" },\n    @{ Old = \"839\u00f78=\"; New = \"182\u00f73=\" },\n    @{ Old = \"568\u00f75=\"; New = \"720\u00f76=\" },\n    @{ Old = \"779\u00f77=\"; New = \"212\u00f79=\" },\n    @{ Old = \"230\u00f75=\"; New = \"837\u00f72=\" },\n    @{ Old = \"356\u00f76=\"; New = \"561\u00f73=\" },\n    @{ Old = \"771\u00f72=\"; New = \"949\u00f77=\" },\n    @{ Old = \"537\u00f77=\"; New = \"953\u00f79=\" }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n}\n"}
